# EIA table_6_01_a.xlsx monthly refresh (2017-01-31 update / energy.gov chunk 7)
# - Header subtitle: "2008 - October 2016" -> "2008 - November 2016"
# - Insert a new "November" data row (row 50) before the closing footnote
#   row, pushing the footnote row from 50 -> 51 and growing the merged
#   footnote range + sheet dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the report period subtitle on row 2.
$ws.Range("A2").Value = "2008 - November 2016"

# 2) Insert a new row above the footnote row (currently row 50), shifting
#    the footnote row (and its merged range) down to row 51.
$ws.Rows("50:50").Insert()

# Copy the formatting of the prior data row (row 49) into the freshly
# inserted row so the new row keeps the same styles (right-aligned number
# format, borders, etc.) as the rest of the monthly data rows.
$ws.Range("A49:D49").Copy()
$ws.Range("A50:D50").PasteSpecial(-4122)

# 3) Populate the new "November" row with its values.
$ws.Range("A50").Value = "November"
$ws.Range("B50").Value = 16636.900000000001
$ws.Range("C50").Value = 12863.5
$ws.Range("D50").Value = 29500.400000000001
